$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Correct column A header name: "Requester_name" -> "Requester_Name"
$ws.Range("A1").Value = "Requester_Name"
